$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update cell values
$ws.Range("E2").Value = 8847544113
$ws.Range("N2").Value = 30990133

# Update the active selection to A2
$ws.Range("A2").Select()
